# Update the "Metadata" worksheet of the ValueSet-presence-valueset workbook
# to reflect the new IG publication metadata (version bump, status change,
# new date, corrected Contact info, and a new Jurisdiction row).

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Metadata")

# --- simple in-place value updates -----------------------------------
$ws.Range("B3").Value = "0.1.7"                                 # Version
$ws.Range("B6").Value = "draft"                                 # Status
$ws.Range("B8").Value = "2024-11-22T12:33:30-06:00"              # Date

# --- shift the tail of the table down by one row to make room for the
#     new "Jurisdiction" row, working from the bottom up so we never
#     clobber a source cell before it has been copied --------------------

# Row 15 ("Immutable" / "BooleanType[null]") moves down to the brand new
# row 16.  Copy formatting from row 15 first so the new row keeps the
# same bordered/wrapped style without Excel inventing a new cellXf.
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)
$ws.Range("A16").Value = "Immutable"
$ws.Range("B16").Value = "BooleanType[null]"

# Row 14 ("Copyright" / empty) moves down to row 15 (now empty value).
$ws.Range("A15").Value = "Copyright"
$ws.Range("B15").ClearContents()

# Row 13 ("Purpose" / empty) moves down to row 14.
$ws.Range("A14").Value = "Purpose"
$ws.Range("B14").ClearContents()

# Row 12 ("Description" / "Permissible values for presence") moves down
# to row 13.
$ws.Range("A13").Value = "Description"
$ws.Range("B13").Value = "Permissible values for presence"

# Row 12 becomes the new "Jurisdiction" row (value left blank).
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").ClearContents()

# --- Contact rows -------------------------------------------------------
# Row 10 keeps the "Contact" label but its value becomes the full
# publisher contact string; row 11 (previously a duplicate "Contact" row)
# now carries the named contact.
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("A11").Value = "Contact"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"
